# Week 16 stat logging + season sim from Week 17
# Applies appended per-game figures to the running play-by-play strings on
# the YDS and ST sheets, and updates the season-total cells on OFF/DEF/ST/
# TURNS/PEN to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append this week's individual play yardages to the
# running space-separated lists.
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$offRush = $wsYDS.Range("B2").Value2
$wsYDS.Range("B2").Value = $offRush + " 2 7 1 3 2 3 14 5 -2 4 6 4 4 5 10 2 3 1 12 2 2 7 2 3 4 5 5 0 -3 4 4 7"

$offPass = $wsYDS.Range("B3").Value2
$wsYDS.Range("B3").Value = $offPass + " 13 0 4 8 -3 12 2 14 20 8 6 8 2 6 4 24 2 23 5 4 8 5 8"

$defRush = $wsYDS.Range("C2").Value2
$wsYDS.Range("C2").Value = $defRush + " 3 2 2 3 2 1 0 8 3 2 9 3 3 -2 3 0"

$defPass = $wsYDS.Range("C3").Value2
$wsYDS.Range("C3").Value = $defPass + " -1 13 9 15 10 8 17 15 9 18 35 5 7 13 7 19 12 4"

# ---------------------------------------------------------------------
# OFF sheet: season totals through Week 16.
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("B2").Value = 6
$wsOFF.Range("C2").Value = 171
$wsOFF.Range("E2").Value = 16
$wsOFF.Range("F2").Value = 51
$wsOFF.Range("G2").Value = 44
$wsOFF.Range("I2").Value = 10
$wsOFF.Range("J2").Value = 28
$wsOFF.Range("O2").Value = 35
$wsOFF.Range("P2").Value = 17

$wsOFF.Range("B3").Value = 9
$wsOFF.Range("C3").Value = 166
$wsOFF.Range("E3").Value = 25
$wsOFF.Range("F3").Value = 91
$wsOFF.Range("G3").Value = 45
$wsOFF.Range("H3").Value = 22
$wsOFF.Range("I3").Value = 55
$wsOFF.Range("J3").Value = 44
$wsOFF.Range("L3").Value = 266
$wsOFF.Range("M3").Value = 171
$wsOFF.Range("Q3").Value = 505

# ---------------------------------------------------------------------
# DEF sheet: season totals through Week 16.
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 210
$wsDEF.Range("F2").Value = 66
$wsDEF.Range("G2").Value = 57
$wsDEF.Range("J2").Value = 26
$wsDEF.Range("N2").Value = 12
$wsDEF.Range("O2").Value = 10
$wsDEF.Range("P2").Value = 6

$wsDEF.Range("C3").Value = 137
$wsDEF.Range("E3").Value = 24
$wsDEF.Range("F3").Value = 101
$wsDEF.Range("G3").Value = 34
$wsDEF.Range("H3").Value = 29
$wsDEF.Range("I3").Value = 51
$wsDEF.Range("J3").Value = 39
$wsDEF.Range("L3").Value = 246
$wsDEF.Range("M3").Value = 166
$wsDEF.Range("Q3").Value = 491

# ---------------------------------------------------------------------
# ST sheet: season totals plus appended per-game kicking lists.
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 55
$wsST.Range("D2").Value = 59
$wsST.Range("F2").Value = 127
$wsST.Range("G2").Value = 119
$wsST.Range("J2").Value = 59
$wsST.Range("K2").Value = 55

$wsST.Range("B3").Value = 26

$koDist = $wsST.Range("B4").Value2
$wsST.Range("B4").Value = $koDist + " 69 66 68"

$koRet = $wsST.Range("B5").Value2
$wsST.Range("B5").Value = $koRet + " 24 23 32"

$ptDist = $wsST.Range("B6").Value2
$wsST.Range("B6").Value = $ptDist + " 25"

$fg40 = $wsST.Range("D3").Value2
$wsST.Range("D3").Value = $fg40 + " 52 57"

$fg50 = $wsST.Range("D4").Value2
$wsST.Range("D4").Value = $fg50 + " 13 0"

$fg60 = $wsST.Range("D5").Value2
$wsST.Range("D5").Value = $fg60 + " 19 0"

# ---------------------------------------------------------------------
# TURNS sheet: season totals through Week 16.
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value = 9
$wsTURNS.Range("D3").Value = 7
$wsTURNS.Range("E3").Value = 7

# ---------------------------------------------------------------------
# PEN sheet: season totals through Week 16.
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 24
